$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("J1").Value = "平均"
$ws.Range("K1").Value = "分散"

# ---- Setosa block (rows 2-5): stats over A2:D51 ----
$ws.Range("H2").Value = "セトナ"
$ws.Range("I2").Value = "がく片"
$ws.Range("J2").Formula = "=AVERAGE(A2:A51)"
$ws.Range("K2").Formula = "=_xlfn.VAR.P(A2:A51)"

$ws.Range("I3").Value = "がく片幅"
$ws.Range("J3").Formula = "=AVERAGE(B2:B51)"
$ws.Range("K3").Formula = "=_xlfn.VAR.P(B2:B51)"

$ws.Range("I4").Value = "花びら"
$ws.Range("J4").Formula = "=AVERAGE(C2:C51)"
$ws.Range("K4").Formula = "=_xlfn.VAR.P(C2:C51)"

$ws.Range("I5").Value = "花びら幅"
$ws.Range("J5").Formula = "=AVERAGE(D2:D51)"
$ws.Range("K5").Formula = "=_xlfn.VAR.P(D2:D51)"

# ---- Versicolor block (rows 7-10): stats over A52:D101 ----
$ws.Range("H7").Value = "バーシクル"
$ws.Range("I7").Value = "がく片"
$ws.Range("J7").Formula = "=AVERAGE(A52:A101)"
$ws.Range("K7").Formula = "=_xlfn.VAR.P(A52:A101)"

$ws.Range("I8").Value = "がく片幅"
$ws.Range("J8").Formula = "=AVERAGE(B52:B101)"
$ws.Range("K8").Formula = "=_xlfn.VAR.P(B52:B101)"

$ws.Range("I9").Value = "花びら"
$ws.Range("J9").Formula = "=AVERAGE(C52:C101)"
$ws.Range("K9").Formula = "=_xlfn.VAR.P(C52:C101)"

$ws.Range("I10").Value = "花びら幅"
$ws.Range("J10").Formula = "=AVERAGE(D52:D101)"
$ws.Range("K10").Formula = "=_xlfn.VAR.P(D52:D101)"

# ---- Virginica block (rows 12-15): stats over A102:D151 ----
$ws.Range("H12").Value = "バージニカ"
$ws.Range("I12").Value = "がく片"
$ws.Range("J12").Formula = "=AVERAGE(A102:A151)"
$ws.Range("K12").Formula = "=_xlfn.VAR.P(A102:A151)"

$ws.Range("I13").Value = "がく片幅"
$ws.Range("J13").Formula = "=AVERAGE(B102:B151)"
$ws.Range("K13").Formula = "=_xlfn.VAR.P(B102:B151)"

$ws.Range("I14").Value = "花びら"
$ws.Range("J14").Formula = "=AVERAGE(C102:C151)"
$ws.Range("K14").Formula = "=_xlfn.VAR.P(C102:C151)"

$ws.Range("I15").Value = "花びら幅"
$ws.Range("J15").Formula = "=AVERAGE(D102:D151)"
$ws.Range("K15").Formula = "=_xlfn.VAR.P(D102:D151)"

# ---- Row 17: note ----
$ws.Range("H17").Value = "2列目のデータを使用"

# ---- Row 18: column labels for the Bayes table ----
$ws.Range("I18").Value = "セトナ"
$ws.Range("J18").Value = "バーシクル"
$ws.Range("K18").Value = "バージニア"

# ---- Row 19: Bayes formulas using A2 (row-1 sample) ----
$ws.Range("H19").Value = "ベイズ式"
$ws.Range("I19").Formula = "=1*(EXP(-(A2-J2)^2)/2*K2)/SQRT(2*PI()*J2)"
$ws.Range("J19").Formula = "=1*(EXP(-(A2-J7)^2)/2*K7)/SQRT(2*PI()*K7)"
$ws.Range("K19").Formula = "=1*(EXP(-(A2-J12)^2)/2*K12)/SQRT(2*PI()*K12)"
$ws.Range("L19").Value = "がく片"

# ---- Row 20 ----
$ws.Range("I20").Formula = "=1*(EXP(-(B2-J3)^2)/2*K3)/SQRT(2*PI()*J3)"
$ws.Range("J20").Formula = "=1*(EXP(-(B2-J8)^2)/2*K8)/SQRT(2*PI()*K8)"
$ws.Range("K20").Formula = "=1*(EXP(-(B2-J13)^2)/2*K13)/SQRT(2*PI()*K13)"
$ws.Range("L20").Value = "がく片幅"

# ---- Row 21 ----
$ws.Range("I21").Formula = "=1*(EXP(-(C2-J4)^2)/2*K4)/SQRT(2*PI()*J4)"
$ws.Range("J21").Formula = "=1*(EXP(-(C2-J9)^2)/2*K9)/SQRT(2*PI()*K9)"
$ws.Range("K21").Formula = "=1*(EXP(-(C2-J14)^2)/2*K14)/SQRT(2*PI()*K14)"
$ws.Range("L21").Value = "花びら"

# ---- Row 22 ----
$ws.Range("I22").Formula = "=1*(EXP(-(D2-J5)^2)/2*K5)/SQRT(2*PI()*J5)"
$ws.Range("J22").Formula = "=1*(EXP(-(D2-J10)^2)/2*K10)/SQRT(2*PI()*K10)"
$ws.Range("K22").Formula = "=1*(EXP(-(D2-J15)^2)/2*K15)/SQRT(2*PI()*K15)"
$ws.Range("L22").Value = "花びら幅"

# ---- Column widths to roughly match the target layout ----
# (ColumnWidth is stored internally as characters + 5/MDW padding, quantized
# to 1/MDW pixels with MDW=7, so these inputs are chosen to land as close as
# possible to the target stored widths of 9.83203125 / 13 / 10.5 / 10.)
$ws.Columns.Item(7).ColumnWidth = 9.142857142857144
$ws.Columns.Item(8).ColumnWidth = 9.142857142857144
$ws.Columns.Item(9).ColumnWidth = 12.285714285714286
$ws.Columns.Item(10).ColumnWidth = 9.857142857142858
$ws.Columns.Item(11).ColumnWidth = 9.285714285714286

# ---- Selection, matching the saved view in the target ----
$ws.Range("I19").Select()
